# Test animations have been implemented
# Insert a new task row at row 72 ("Popup Menu Bug") on the milestones sheet,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the current row 72 (KI System), shifting everything down.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row with the new task.
$ws.Range("A72").Value = "Popup Menu Bug"
$ws.Range("B72").Value = "Fade Popup out when specific object isn't active"
$ws.Range("C72").Value = "Graphics"
$ws.Range("D72").Value = "very high"
$ws.Range("E72").Value = "my-reality@gmx.de"
$ws.Range("G72").Formula = "=ISNUMBER(F72)"

# Keep the view roughly where the author left it.
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Range("H72").Select()
